# Commit: "Don't need to convert with fuhep"
# The N column ("predicted Clint" style calc) previously multiplied by I
# (the fu,hep fraction-unbound-in-hepatocyte-prep column) before scaling;
# that multiplication is removed so the formula is just J/120/40*1000.
# The (previously blank) rows where J was empty now also get the formula
# (evaluating to 0), and the whole N2:N58 range picks up a "0.00" number
# format that didn't exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N2 keeps a standalone (non-shared) formula.
$ws.Range("N2").Formula = "=J2/120/40*1000"

# N3:N58 becomes one filled-down block (formerly only N4:N58 had the
# formula; N3 and several other rows with blank J cells were skipped).
$ws.Range("N3:N58").Formula = "=J3/120/40*1000"

# Apply the new "0.00" number format (a fresh cellXfs entry) across the
# whole column range that now carries values.
$ws.Range("N2:N58").NumberFormat = "0.00"

# Update the active selection / scroll position to match the saved view.
$ws.Range("O8").Select()
